$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value to literal text without triggering Excel's
# automatic number/percent parsing and without creating a new cell style.
# We do this by writing the text into a scratch cell via a formula (so it
# becomes a plain calculated string), then copy/paste-special (values only)
# that computed value onto the destination - this preserves the
# destination's existing style while landing a literal text value.
function Set-LiteralText {
    param(
        [string]$Address,
        [string]$Text
    )
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '="' + $Text + '"'
    $scratch.Copy()
    $ws.Range($Address).PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
    $excel.CutCopyMode = $false
}

# ---- Row 2: reorder "Recorded By" list ----
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# ---- Row 3: reorder "Recorded By" list ----
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# ---- Row 4: reorder "Recorded By" list ----
$ws.Range("G4").Value = "servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# ---- Row 5: reorder "Recorded By" list ----
$ws.Range("G5").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# ---- Row 6: Recorded Sessions metric 12 -> 13 ----
$ws.Range("L6").Value = 13

# ---- Row 7: reorder "Recorded By" list ----
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"

# ---- Row 7: Missing Sessions metric 2 -> 1 ----
$ws.Range("L7").Value = 1

# ---- Row 9: reorder "Recorded By" list ----
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# ---- Row 9: Coverage % 41.4% -> 44.8% (keep as literal text, same style) ----
Set-LiteralText "L9" "44.8%"

# ---- Row 10: Average Attendance % 24.1% -> 24.5% (keep as literal text, same style) ----
Set-LiteralText "L10" "24.5%"

# ---- Row 12: reorder "Recorded By" list ----
$ws.Range("G12").Value = "yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg"

# ---- Row 15: reorder "Recorded By" list ----
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# ---- Row 15: summary metrics ----
$ws.Range("O15").Value = 13
$ws.Range("P15").Value = 1
Set-LiteralText "R15" "44.8%"
Set-LiteralText "S15" "24.5%"

# ---- Row 27: session became recorded; copy row 28's formatting (style "2")
#      onto row 27 (previously styled as "Not Recorded" / style "9") ----
$ws.Range("A28:I28").Copy()
$ws.Range("A27:I27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg"
$ws.Range("H27").Value = "76/251"
$ws.Range("I27").Value = "Recorded"

# ---- Row 28: reorder "Recorded By" list ----
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
